$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = '1327341'
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1327341'
$ws.Range("C2").Value = '[Impact Brazil] - Content Creator & Social Media'
$ws.Range("D2").Value = 'São Paulo, SP, Brasil'
$ws.Range("F2").Value = '2 applicants'
$ws.Range("H2").Value = 'Elebbre'

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = '1327303'
$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1327303'
$ws.Range("C3").Value = 'Project Coordinator Intern'
$ws.Range("D3").Value = 'Chandigarh, India'
$ws.Range("G3").Value = '3 - 6 Months'

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = '1327205'
$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1327205'
$ws.Range("C4").Value = 'Environmental Technician'
$ws.Range("D4").Value = 'Calgary, AB, Canada'
$ws.Range("G4").Value = '6 - 18 Months'
$ws.Range("H4").Value = 'Oak Environmental Inc.'

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = '1327139'
$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1327139'
$ws.Range("C5").Value = 'Desktop Software Engineer'
$ws.Range("D5").Value = 'Ahmedabad, Gujarat, India'
$ws.Range("H5").Value = 'WeHear Innovations Private Limited'

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = '1327128'
$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1327128'
$ws.Range("C6").Value = 'Digital Marketing'
$ws.Range("D6").Value = 'Καλλιθέα 630 77, Ελλάδα'
$ws.Range("F6").Value = '5 applicants'
$ws.Range("H6").Value = 'Respirotours'

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '1326643'
$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1326643'
$ws.Range("C7").Value = 'On Premise Analyst'
$ws.Range("D7").Value = 'Panamá, Provincia de Panamá, Panamá'
$ws.Range("F7").Value = '98 applicants'
$ws.Range("G7").Value = '6 - 18 Months'
$ws.Range("H7").Value = 'Red Bull Panamá'

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '1326467'
$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1326467'
$ws.Range("C8").Value = 'Customer Relations Executive/ Intern'
$ws.Range("D8").Value = 'Colombo, Sri Lanka'
$ws.Range("F8").Value = '65 applicants'
$ws.Range("G8").Value = '6 - 18 Months'
$ws.Range("H8").Value = 'Formix (Pvt) Ltd'

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '1325417'
$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1325417'
$ws.Range("C9").Value = 'Junior Software Engineer – AI & Internal Tools (EU ONLY)'
$ws.Range("D9").Value = 'Brussels, Belgium'
$ws.Range("F9").Value = '89 applicants'
$ws.Range("G9").Value = '6 - 18 Months'
$ws.Range("H9").Value = 'Eureka Resource Mining'

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '1325130'
$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1325130'
$ws.Range("C10").Value = 'Marketing specialist'
$ws.Range("D10").Value = 'Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt'
$ws.Range("F10").Value = '11 applicants'
$ws.Range("H10").Value = 'print shop'

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '1324500'
$ws.Range("B11").Value = 'https://aiesec.org/opportunity/global-talent/1324500'
$ws.Range("C11").Value = 'DT Software Engineer Trainee (EU ONLY)'
$ws.Range("D11").Value = 'Brussels, Belgium'
$ws.Range("H11").Value = 'UCB'

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '1322343'
$ws.Range("B12").Value = 'https://aiesec.org/opportunity/global-talent/1322343'
$ws.Range("C12").Value = 'Electronics and Electrical Intern'
$ws.Range("D12").Value = 'Manipal, Karnataka, India'
$ws.Range("F12").Value = '21 applicants'
$ws.Range("G12").Value = '3 - 6 Months'
$ws.Range("H12").Value = 'M.A.H.E.'

# Row 13 (new row)
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = '1310154'
$ws.Range("B13").Value = 'https://aiesec.org/opportunity/global-talent/1310154'
$ws.Range("C13").Value = 'UX/UI Designer'
$ws.Range("D13").Value = 'Mississauga, Canada'
$ws.Range("E13").Value = 'No'
$ws.Range("F13").Value = '240 applicants'
$ws.Range("G13").Value = '6 - 18 Months'
$ws.Range("H13").Value = 'Remitbee'

# Row 14 (new row)
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = '1304736'
$ws.Range("B14").Value = 'https://aiesec.org/opportunity/global-talent/1304736'
$ws.Range("C14").Value = 'Guest Relations Manager'
$ws.Range("D14").Value = 'Heraklion, Greece'
$ws.Range("E14").Value = 'No'
$ws.Range("F14").Value = '177 applicants'
$ws.Range("G14").Value = '9 - 12 Weeks'
$ws.Range("H14").Value = 'Remarc Internation'

# Column width changes
$ws.Columns("C").ColumnWidth = 58.17
$ws.Columns("D").ColumnWidth = 56.17
$ws.Columns("H").ColumnWidth = 36.17
